# Update automàtic: dades i banners [2026-02-20 09:15]
#
# This script reproduces, on the "Dades_Període" worksheet, the data refresh
# captured by the XML diff: new measurement period (08:30 - 09:00), refreshed
# extraction timestamps, and updated weather readings for the first data row
# (row 2), plus refreshed extraction timestamps on rows 3-6.
#
# Numeric-looking values are written with a leading apostrophe so Excel keeps
# them stored as text (matching the source workbook, where every data cell is
# an inline/shared string rather than a numeric cell) instead of silently
# converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Període")

# --- Row 2 -----------------------------------------------------------------
$ws.Range("E2").Value  = "08:30 - 09:00"
$ws.Range("H2").Value  = "2026-02-20 09:15:32"

$ws.Range("M2").Value  = "'272"
$ws.Range("N2").Value  = "'64"
$ws.Range("O2").Value  = "'1023.9"
$ws.Range("P2").Value  = "'0.1"
$ws.Range("Q2").Value  = "08:30 - 09:00"
$ws.Range("R2").Value  = "'381"
$ws.Range("S2").Value  = "'9.2"
$ws.Range("T2").Value  = "'8.6"
$ws.Range("U2").Value  = "'10.3"
$ws.Range("V2").Value  = "'1.1"
$ws.Range("W2").Value  = "'4.0"

$ws.Range("X2").Value  = "08:30 - 09:00"
$ws.Range("Y2").Value  = "'9.2"
$ws.Range("Z2").Value  = "'10.3"
$ws.Range("AA2").Value = "'8.6"
$ws.Range("AB2").Value = "'64"
$ws.Range("AC2").Value = "'0.1"
$ws.Range("AD2").Value = "'1.1"
$ws.Range("AE2").Value = "'272"
$ws.Range("AF2").Value = "'4.0"
$ws.Range("AG2").Value = "'1023.9"
$ws.Range("AH2").Value = "'381"

$ws.Range("AI2").Value = "08:30 - 09:00"
$ws.Range("AJ2").Value = "'9.2"
$ws.Range("AK2").Value = "'10.3"
$ws.Range("AL2").Value = "'8.6"
$ws.Range("AM2").Value = "'64"
$ws.Range("AN2").Value = "'0.1"
$ws.Range("AO2").Value = "'1.1"
$ws.Range("AP2").Value = "'272"
$ws.Range("AQ2").Value = "'4.0"
$ws.Range("AR2").Value = "'1023.9"
$ws.Range("AS2").Value = "'381"

# --- Rows 3-6: refreshed extraction timestamp only --------------------------
$ws.Range("H3").Value = "2026-02-20 09:15:34"
$ws.Range("H4").Value = "2026-02-20 09:15:34"
$ws.Range("H5").Value = "2026-02-20 09:15:34"
$ws.Range("H6").Value = "2026-02-20 09:15:34"
